$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting of the last existing data row (229, column A)
# onto the four new date cells before writing their values.
$ws.Range("A229").Copy() | Out-Null
$ws.Range("A230:A233").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row = 230; Date = 44304; B = 1; C = 1; D = 62.34413965087282 },
    @{ Row = 231; Date = 44305; B = 0; C = 1; D = 62.34413965087282 },
    @{ Row = 232; Date = 44306; B = 0; C = 1; D = 62.34413965087282 },
    @{ Row = 233; Date = 44307; B = 0; C = 1; D = 62.34413965087282 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
